# Apply the "new version with timestamp" update to the DaySale report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. ATROVENT 500MCG/2ML 20 UNIT DOSE VIALS (row 14): ratio 1:5 -> 1:6
$ws.Range("H14").Value = "1:6"

# 2. DOLIPRANE 1 GM 15 TABS. (row 29): ratio 9:1 -> 8:2, price 31.6800 -> 63.8400,
#    and number-of-transactions ratio 0:2 -> 1:1
# (P29 is a numeric-formatted cell that actually stores its price as text, so a
#  leading apostrophe is used to force a text value, same as typing it in Excel.)
$ws.Range("H29").Value = "8:2"
$ws.Range("P29").Value = "'63.8400"
$ws.Range("Q29").Value = "1:1"

# 3. PULMICORT 0.5MG/ML 20 NEBULIZER VIAL SUSP. (row 62): ratio 1:15 -> 1:16
$ws.Range("H62").Value = "1:16"

# 4. Remove the "VIOTIC EAR DROPS 10 ML" line (row 78) completely; this shifts
#    every row below it up by one, exactly like a manual row delete in Excel.
$ws.Rows(78).Delete()

# 5. The item index numbers in column A for the rows that moved up (formerly
#    rows 79-94, now rows 78-93) must be decremented by one since an item was
#    removed from the list.
for ($r = 78; $r -le 93; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = [double]$current - 1
}

# 6. Recalculate the total shown at the bottom of the "sell price" column
#    (now row 94 after the row shift above).
$total = 0.0
for ($r = 7; $r -le 93; $r++) {
    $cell = $ws.Cells.Item($r, 16)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $total = $total + [double]$val
    }
}
$ws.Range("P94").Value = $total

# 7. Update the generated timestamp shown in the footer (now row 95 after the
#    row shift above) to reflect the new export time.
$ws.Range("A95").Value = "Wednesday, 3 September, 2025 10:00 PM"
